# regen sval data to filter save games
# Update the stat columns (B:G) for rows 2-15 on the active sheet with the
# regenerated values (save games are now filtered out of the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126),
    @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 1, 10.35301142835362),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538),
    @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 3.811642989160245),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(0.3464964993005633, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 0, 1.611428759079651),
    @(0.02258322285507441, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 2.328612170846122),
    @(0.1554434735375247, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 3.021299490525381),
    @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 4.371470058157054),
    @(0.06328177979961902, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 0, 1.613486173897131),
    @(0.7287194209349384, 1.65323645889881, 157.8057217802531, 6.48142807727062, 1, 166.6691057373575),
    @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 6.741336633845642),
    @(0.001754667048134761, 0.004309184025731883, 0.1529057820181812, 0.4998867070740569, 1, 0.6588563401661047)
)

$startRow = 2
$cols = @("B", "C", "D", "E", "F", "G")
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$r").Value = $vals[$j]
    }
}
